# Weekly data refresh for "Hortaliza, Vega Monumental Concepción - Ajo":
# a new week's record is inserted as row 97, and every subsequent record
# (old rows 97..163) shifts down by one row (new rows 98..164). The last
# existing record (old row 163) ends up duplicated into the newly
# appended row 164.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFmt = "YYYY-MM-DD HH:MM:SS"

# Columns that carry record-specific data and therefore shift down.
# (A, B, C, E, F, G, H, R hold the constant "market/category" values
# shared by every row in this subset and are left untouched.)
$shiftCols = 4, 9, 10, 11, 12, 13, 14, 15, 16, 17   # D, I, J, K, L, M, N, O, P, Q

# Walk from the bottom (164) up to 98 so each source row is read with
# Value2 before it gets overwritten by the row above it shifting down.
for ($destRow = 164; $destRow -ge 98; $destRow--) {
    $srcRow = $destRow - 1
    foreach ($col in $shiftCols) {
        $val = $ws.Cells.Item($srcRow, $col).Value2
        if ($col -eq 4) {
            $ws.Cells.Item($destRow, $col).NumberFormat = $dateFmt
        }
        $ws.Cells.Item($destRow, $col).Value = $val
    }
}

# New row 164 also needs the constant columns populated (same values as
# every other row in this subset).
$ws.Cells.Item(164, 1).Value = 11
$ws.Cells.Item(164, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(164, 3).Value = "Bíobío"
$ws.Cells.Item(164, 5).Value = 8
$ws.Cells.Item(164, 6).Value = 100112003
$ws.Cells.Item(164, 7).Value = "Ajo"
$ws.Cells.Item(164, 8).Value = "Chino"
$ws.Cells.Item(164, 18).Value = "Hortaliza"

# Finally, row 97 receives the brand-new weekly record.
$ws.Cells.Item(97, 4).NumberFormat = $dateFmt
$ws.Cells.Item(97, 4).Value = 44673
$ws.Cells.Item(97, 10).Value = 220
$ws.Cells.Item(97, 11).Value = 18000
$ws.Cells.Item(97, 12).Value = 19000
$ws.Cells.Item(97, 13).Value = 18455
$ws.Cells.Item(97, 16).Value = 1846
